# Generate Report for Handback
# The second file (a9abe794-abb3-4b49-bd0d-95c52ceed991.md) has now been
# handed back and is in sync with en-US, so its status changes from
# "Ready for handoff" to "Handed back: in sync with en-US" on every
# sheet, and the per-locale "Latest Handback DateTime" is refreshed.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $status
$wsZh.Range("H2").Value = "2016-03-19 17:22:48"
$wsZh.Range("H3").Value = "2016-03-19 17:22:48"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $status
$wsDe.Range("H2").Value = "2016-03-19 17:23:01"
$wsDe.Range("H3").Value = "2016-03-19 17:23:01"
